$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 2

$replacements = @(
  @{ Old = 'Singathanda ukuba nencoko ngomxeba nawe ozakuthatha imizuzu engamashumi amane anesihlanu ubude. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. Ukuba ngelixa wenziwa udliwano-ndlebe, uye waphazamiseka, ndicela ucele umphandi ukuba ame ude uzive ukhuselekile ukuqhubeka nokuthetha.'; New = 'Singathanda ukuba nencoko ngomxeba nawe ozakuthatha imizuzu engamashumi amane anesihlanu ubude. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Nceda uqinisekise ukuba xa sikutsalela umnxeba, uvumela kuphela udliwano-ndlebe ukuba luqale xa ukwindawo yabucala apho uziva ukhululekile ukuthetha ngaphandle kokumanyelwa okanye ukuphazanyiswa. Ukuba ngelixa wenziwa udliwano-ndlebe, uye waphazamiseka, ndicela ucele umphandi ukuba ame ude uzive ukhuselekile ukuqhubeka nokuthetha.' },
  @{ Old = 'We would like to know more about your experience with the Keeping My Child Safe module. Your opinions will help us improve the module, and hopefully improve the experience of other parents like you in the future. '; New = 'Singathanda ukwazi ngakumbi malunga namava akho ngemodyuli yokuGcina Umntwana Wam eKhuselekile. Izimvo zakho ziya kusinceda siphucule imodyuli, kwaye sinethemba lokuphucula amava abanye abazali abafana nawe kwixesha elizayo. ' },
  @{ Old = 'Do I have to agree to be interviewed?'; New = 'Ingaba kufuneka ndivume ukuba nodliwano-ndlebe?' },
  @{ Old = 'What happens with my information?'; New = 'Kwenzeka ntoni ngenkcukacha zam?' },
  @{ Old = 'To protect your personal information (including your real name, contact details, and any other information that can identify you), we will give you a participant number, and you can choose a name you want us to call you during the interview.'; New = 'Ukukhusela iinkcukacha zakho (kuquka igama lakho lokwenene, inkcukacha zoqhagamshelwano, kunye naluphi na olunye ulwazi olungachaza wena), sizakunika inombolo yokuthatha inxaxheba, kwaye ungazikhethela igama ofuna sikubize ngalo ngexesha lodliwano-ndlebe.' },
  @{ Old = 'With your permission, we will record the interview to help us remember the discussion and later write down what was said. We will delete any personal information we collect from you at the end of the study and, after transcribing your interview, change any data which might lead to identification at the point of transcription. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijonge oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kuqwalaselwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu eqela lophando kuphela agunyazisiweyo aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa.'; New = 'Ngemvume yakho, siya kurekhoda udliwano-ndlebe ukuze usincede sikhumbule ingxoxo kwaye kamva sibhale phantsi oko bekuthethiwe. Siya kucima nayiphi na ingcaciso yobuqu esiyiqokelele kuwe ekupheleni kophononongo kwaye, emva kokubhala udliwano-ndlebe lwakho, sitshintshe nayiphi na idatha enokukhokelela ekuchongeni kwindawo yokukhuphela. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijonge oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kuqwalaselwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu eqela lophando kuphela agunyazisiweyo aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa.' },
  @{ Old = 'As a thank you for taking part in the discussion, we''ll give you a R30 airtime voucher/data bundle. '; New = 'Njengombulelo ngokuthatha inxaxheba kule ngxoxo, siza kukunika ivawutsha ye-airtime/data bundle ye-R30. ' },
  @{ Old = 'What happens to my information if I agree to be interviewed?'; New = 'Kwenzeka ntoni kwinkcukacha zam ukuba ndiyavuma ukuba noludliwano-ndlebe?' },
  @{ Old = 'Siqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcina ngokukhuselekileyo. Your information, like your consent form and interview recording, and any information you provide via email or WhatsApp, will be kept safe on secure servers at the University of Cape Town. '; New = 'Siqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcina ngokukhuselekileyo. Ulwazi lwakho, olufana nefomu yakho yemvume kunye nokurekhodwa kodliwano-ndlebe, nayo nayiphi na ingcaciso oyinikeza nge-imeyile okanye nge-WhatsApp, iya kugcinwa ikhuselekile kwiiseva ezikhuselekileyo kwiDyunivesithi yaseKapa. ' },
  @{ Old = 'Interview recordings will be deleted after we have written our notes. Any details that identify you will be kept separate and only authorised staff can access them. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kususwa xa isifundo siphelile. '; New = 'Interview recordings will be deleted after we have written our notes. Naziphi na iinkcukacha ezichongayo ziya kugcinwa ngokwahlukileyo kwaye ngabasebenzi abagunyazisiweyo kuphela abanokufikelela kuzo. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kususwa xa isifundo siphelile. ' },
  @{ Old = 'Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. We plan to share the results in reports and at conferences so others can learn from this study too.'; New = 'Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. Sicwangcisa ngokwabelana ngeziphumo kwiingxelo nakwii-nkomfa ukuze nabanye bafunde kolu phononongo.' },
  @{ Old = 'The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town.'; New = 'Abaphandi abaziintloko kolu phononongo nguProf Cathy Ward no Cindee Bruyns ze uCo-investigator ibengu Carly Katzef bonke basuka kwiDyunivesithi yaseKapa.' },
  @{ Old = 'Who pays for the study?'; New = 'Ngubani obhatalela oluphononongo?' },
  @{ Old = 'This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. '; New = 'Olu phononongo luyinxalenye ye Global Parenting Initiative, luxhaswe ngokwezimali yi LEGO Foundation, i-Oak Fundation, i-World Childhood Foundation, i-Human Safety Net kunye ne UK Research kunye ne Innovaion Global Challenges Research Fund. ' },
  @{ Old = 'The University Cape Town makes sure your personal information is used safely and correctly, just for research. Uphononongo lulandela imithetho yokukhuselwa kwedatha efana ne-GDPR (General Data Protection Regulation) e-UK kunye ne-POPIA (uMthetho woKhuselo loLwazi loMntu) eMzantsi Afrika. Nayiphi na idatha ethi ithunyelwe ngaphesheya kwemida izakuthobelana ne POPIA. '; New = 'I Dyunivesithi yaseKapa iqinisekisa ukuba iinkcukacha zakho zobuqu zisetyenziswa ngokukhuselekileyo nangokuchanekileyo, nje kuphando kuphela. Uphononongo lulandela imithetho yokukhuselwa kwedatha efana ne-GDPR (General Data Protection Regulation) e-UK kunye ne-POPIA (uMthetho woKhuselo loLwazi loMntu) eMzantsi Afrika. Nayiphi na idatha ethi ithunyelwe ngaphesheya kwemida izakuthobelana ne POPIA. ' },
  @{ Old = '[Yakube ivunyiwe imigaqo yokuziphatha iya kufundeka ngoluhlobo lulandelayo: Olu phononongo lufumene imvume kwiDyunivesithi yaseKapa kwiZiko leKomiti yeeNqobo zoPhando kwiNzululwazi yezeNtlalo kunye neKomiti yeeNqobo zokuziphatha zoPhando lweDyunivesithi yaseKapa. The study has also been approved by the Western Cape Department of Health and Wellness and Department of Social Development, and City of Cape Town’s City health.]'; New = '[Yakube ivunyiwe imigaqo yokuziphatha iya kufundeka ngoluhlobo lulandelayo: Olu phononongo lufumene imvume kwiDyunivesithi yaseKapa kwiZiko leKomiti yeeNqobo zoPhando kwiNzululwazi yezeNtlalo kunye neKomiti yeeNqobo zokuziphatha zoPhando lweDyunivesithi yaseKapa. Olu phononongo lukwavunyiwe liSebe lezeMpilo neMpilo eNtshona Koloni kunye neSebe loPhuhliso loLuntu, kunye nezempilo kwiSixeko saseKapa.]' },
  @{ Old = 'Please respond with the word “agree” to each as I go through each of the following points. If you don’t agree, we can go over any other information you need to make your decision and if you still agree then we can proceed:'; New = 'Nceda uphendule ngegama elithi "ndiyavuma" kumntu ngamnye njengoko ndihamba ngenqaku ngalinye kula alandelayo. Ukuba awuvumelani, singadlula kulo naliphi na ulwazi oludingayo ukuze uthathe isigqibo kwaye ukuba usavuma singaqhubeka ke:' },
)

$idx = 0
foreach ($r in $replacements) {
  $idx = $idx + 1
  $rng = $d.Content
  $found = $rng.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $r.New, $wdReplaceOne)
  if (-not $found) {
    Write-Host "NOT FOUND #${idx}: $($r.Old)"
  } else {
    Write-Host "Replaced #${idx}"
  }
}
